# Updated prices in master parts list. Closes #160.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Used in Sections" (column K) text re-orderings ---------------------

$ws.Range("K133").Value = "Wheel Assembly, Rocker-Bogie, Differential Pivot, Mechanical Integration, Body, Corner Steering"
$ws.Range("K134").Value = "Wheel Assembly, Differential Pivot"
$ws.Range("K137").Value = "Wheel Assembly, Rocker-Bogie, Differential Pivot, Head Assembly"
$ws.Range("K139").Value = "Rocker-Bogie, Wheel Assembly, Mechanical Integration, Differential Pivot"

# --- Price updates (column I = Price Each, column J = Price Total) -------

$ws.Range("I12").Value = 3.02
$ws.Range("J12").Value = 3.02

$ws.Range("I80").Value = 25.74
$ws.Range("J80").Value = 25.74

$ws.Range("I82").Value = 2.19
$ws.Range("J82").Value = 10.95

$ws.Range("I83").Value = 1.59
$ws.Range("J83").Value = 1.59

$ws.Range("I87").Value = 4.99
$ws.Range("J87").Value = 59.88

$ws.Range("I89").Value = 6.99
$ws.Range("J89").Value = 41.94

$ws.Range("I90").Value = 5.2
$ws.Range("J90").Value = 20.8

$ws.Range("I93").Value = 5.28
$ws.Range("J93").Value = 5.28

$ws.Range("I113").Value = 2.89
$ws.Range("J113").Value = 23.12

$ws.Range("I114").Value = 1.49
$ws.Range("J114").Value = 1.49

$ws.Range("I118").Value = 2.59
$ws.Range("J118").Value = 10.36

$ws.Range("I122").Value = 1.89
$ws.Range("J122").Value = 3.78

$ws.Range("I123").Value = 3.39
$ws.Range("J123").Value = 13.56

$ws.Range("I127").Value = 0.64
$ws.Range("J127").Value = 0.64

$ws.Range("I128").Value = 1
$ws.Range("J128").Value = 4
